# Andre Russell IPL-KKR scraped-data refresh:
#  - inserts two new columns ("ownTeam", "oppTeam") after "result"
#  - re-scrapes totals (runs/balls/4s/6s/sr) for the existing matches
#  - adds three newly scraped matches (rows 8-10)
#  - re-orders the rows chronologically
# Simplest faithful reproduction: rebuild the whole used range from scratch.

$data = @(
    @("venue", "date", "result", "ownTeam", "oppTeam", "batsman", "totalRuns", "totalBalls", "total4s", "total6s", "sr"),
    @(" Abu Dhabi", " October 07 2020", "KKR won by 10 runs", "Kolkata Knight Riders", "Chennai Super Kings", "Andre Russell ", "2", "4", "0", "0", "50.00"),
    @(" Dubai (DSC)", " September 30 2020", "KKR won by 37 runs", "Kolkata Knight Riders", "Rajasthan Royals", "Andre Russell ", "24", "14", "0", "3", "171.42"),
    @(" Abu Dhabi", " October 16 2020", "Mumbai won by 8 wickets (with 19 balls remaining)", "Kolkata Knight Riders", "Mumbai Indians", "Andre Russell ", "12", "9", "1", "1", "133.33"),
    @(" Abu Dhabi", " October 18 2020", "Match tied (KKR won the one-over eliminator)", "Kolkata Knight Riders", "Sunrisers Hyderabad", "Andre Russell ", "9", "11", "1", "0", "81.81"),
    @(" Abu Dhabi", " October 10 2020", "KKR won by 2 runs", "Kolkata Knight Riders", "Kings XI Punjab", "Andre Russell ", "5", "3", "1", "0", "166.66"),
    @(" Dubai (DSC)", " November 01 2020", "KKR won by 60 runs", "Kolkata Knight Riders", "Rajasthan Royals", "Andre Russell ", "25", "11", "1", "3", "227.27"),
    @(" Sharjah", " October 03 2020", "Capitals won by 18 runs", "Kolkata Knight Riders", "Delhi Capitals", "Andre Russell ", "13", "8", "1", "1", "162.50"),
    @(" Abu Dhabi", " September 23 2020", "Mumbai won by 49 runs", "Kolkata Knight Riders", "Mumbai Indians", "Andre Russell ", "11", "11", "2", "0", "100.00"),
    @(" Sharjah", " October 12 2020", "RCB won by 82 runs", "Kolkata Knight Riders", "Royal Challengers Bangalore", "Andre Russell ", "16", "10", "2", "1", "160.00")
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop everything (old 9-column x 7-row table) so no stale cells survive
# outside the new 11-column x 10-row range.
$ws.Cells.Clear()

$numRows = $data.Count
$numCols = $data[0].Count

for ($r = 0; $r -lt $numRows; $r++) {
    for ($c = 0; $c -lt $numCols; $c++) {
        $cell = $ws.Cells.Item($r + 1, $c + 1)
        # Every value in this sheet is scraped text (numbers-as-text,
        # hence the original file's numberStoredAsText ignoredError), so
        # force the Text format before writing to stop numeric-looking
        # strings like "2" or "50.00" from being coerced to numbers.
        $cell.NumberFormat = "@"
        $cell.Value = $data[$r][$c]
    }
}
